$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = "1000 Bs = 2.78 = 10305.56 pesos"
$newLine1 = "1000 Bs = 2.89 = 10592.96 pesos"
$oldLine2 = "10305.56 pesos = 2.8 = 939.3 Bs"
$newLine2 = "10592.96 pesos = 2.86 = 951.0 Bs"

[string]$text = $ws1.Range("A1").Value()
$text = $text -replace [regex]::Escape($oldLine1), $newLine1
$text = $text -replace [regex]::Escape($oldLine2), $newLine2
$ws1.Range("A1").Value = $text

# --- Update the "tasas" sheet numeric values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 345.89
$ws2.Range("O10").Value = 3664
$ws2.Range("N12").Value = 3699.95
$ws2.Range("O12").Value = 332.17
